$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Distributions")
$ws.Rows("4:5").Insert()
Write-Host "done insert"
$v = $ws.Range("F6").Value2
Write-Host "F6 (was F4) value2: $v"
$v2 = $ws.Range("A7").Formula
Write-Host "A7 formula: $v2"
